# Add Sheet2 as a copy of Sheet1 (preserves merges, formatting, dimension),
# then rewrite the data cells to reflect the new (rebel-shifted) values,
# and fix up the view/selection state on both sheets.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Duplicate Sheet1 right after itself -> becomes "Sheet1 (2)"
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Sheet2"

$ws2.Range("B3").Value = "A13"
$ws2.Range("C3").Value = "C13"
$ws2.Range("D3").Value = "D13"
$ws2.Range("E3").Value = "E13"
$ws2.Range("F3").Value = "F13"
$ws2.Range("G3").Value = "G13"
$ws2.Range("H3").Value = "H13"
$ws2.Range("B4").Value = "A14"
$ws2.Range("C4").Value = "C14"
$ws2.Range("D4").Value = "D14"
$ws2.Range("E4").Value = "E14"
$ws2.Range("F4").Value = "F14"
$ws2.Range("G4").Value = "G14"
$ws2.Range("H4").Value = "H14"
$ws2.Range("B5").Value = "A15"
$ws2.Range("C5").Value = "C15"
$ws2.Range("D5").Value = "D15"
$ws2.Range("E5").Value = "E15"
$ws2.Range("F5").Value = "F15"
$ws2.Range("G5").Value = "G15"
$ws2.Range("H5").Value = "H15"
$ws2.Range("B6").Value = "A16"
$ws2.Range("C6").Value = "C16"
$ws2.Range("D6").Value = "D16"
$ws2.Range("E6").Value = "E16"
$ws2.Range("F6").Value = "F16"
$ws2.Range("G6").Value = "G16"
$ws2.Range("H6").Value = "H16"
$ws2.Range("B7").Value = "A17"
$ws2.Range("C7").Value = "C17"
$ws2.Range("D7").Value = "D17"
$ws2.Range("E7").Value = "E17"
$ws2.Range("F7").Value = "F17"
$ws2.Range("G7").Value = "G17"
$ws2.Range("H7").Value = "H17"
$ws2.Range("B8").Value = "A18"
$ws2.Range("C8").Value = "C18"
$ws2.Range("D8").Value = "D18"
$ws2.Range("E8").Value = "E18"
$ws2.Range("F8").Value = "F18"
$ws2.Range("G8").Value = "G18"
$ws2.Range("H8").Value = "H18"
$ws2.Range("B9").Value = "A19"
$ws2.Range("C9").Value = "C19"
$ws2.Range("D9").Value = "D19"
$ws2.Range("E9").Value = "E19"
$ws2.Range("F9").Value = "F19"
$ws2.Range("G9").Value = "G19"
$ws2.Range("H9").Value = "H19"
$ws2.Range("B10").Value = "A20"
$ws2.Range("C10").Value = "C20"
$ws2.Range("D10").Value = "D20"
$ws2.Range("E10").Value = "E20"
$ws2.Range("F10").Value = "F20"
$ws2.Range("G10").Value = "G20"
$ws2.Range("H10").Value = "H20"
$ws2.Range("B11").Value = "A21"
$ws2.Range("C11").Value = "C21"
$ws2.Range("D11").Value = "D21"
$ws2.Range("E11").Value = "E21"
$ws2.Range("F11").Value = "F21"
$ws2.Range("G11").Value = "G21"
$ws2.Range("H11").Value = "H21"
$ws2.Range("B12").Value = "A22"
$ws2.Range("C12").Value = "C22"
$ws2.Range("D12").Value = "D22"
$ws2.Range("E12").Value = "E22"
$ws2.Range("F12").Value = "F22"
$ws2.Range("G12").Value = "G22"
$ws2.Range("H12").Value = "H22"
$ws2.Range("B13").Value = "A23"
$ws2.Range("C13").Value = "C23"
$ws2.Range("D13").Value = "D23"
$ws2.Range("E13").Value = "E23"
$ws2.Range("F13").Value = "F23"
$ws2.Range("G13").Value = "G23"
$ws2.Range("H13").Value = "H23"
$ws2.Range("B14").Value = "A24"
$ws2.Range("C14").Value = "C24"
$ws2.Range("D14").Value = "D24"
$ws2.Range("E14").Value = "E24"
$ws2.Range("F14").Value = "F24"
$ws2.Range("G14").Value = "G24"
$ws2.Range("H14").Value = "H24"
$ws2.Range("B15").Value = "A25"
$ws2.Range("C15").Value = "C25"
$ws2.Range("D15").Value = "D25"
$ws2.Range("E15").Value = "E25"
$ws2.Range("F15").Value = "F25"
$ws2.Range("G15").Value = "G25"
$ws2.Range("H15").Value = "H25"
$ws2.Range("B16").Value = "A26"
$ws2.Range("C16").Value = "C26"
$ws2.Range("D16").Value = "D26"
$ws2.Range("E16").Value = "E26"
$ws2.Range("F16").Value = "F26"
$ws2.Range("G16").Value = "G26"
$ws2.Range("H16").Value = "H26"
$ws2.Range("B17").Value = "A27"
$ws2.Range("C17").Value = "C27"
$ws2.Range("D17").Value = "D27"
$ws2.Range("E17").Value = "E27"
$ws2.Range("F17").Value = "F27"
$ws2.Range("G17").Value = "G27"
$ws2.Range("H17").Value = "H27"
$ws2.Range("B18").Value = "A28"
$ws2.Range("C18").Value = "C28"
$ws2.Range("D18").Value = "D28"
$ws2.Range("E18").Value = "E28"
$ws2.Range("F18").Value = "F28"
$ws2.Range("G18").Value = "G28"
$ws2.Range("H18").Value = "H28"
$ws2.Range("B19").Value = "A29"
$ws2.Range("C19").Value = "C29"
$ws2.Range("D19").Value = "D29"
$ws2.Range("E19").Value = "E29"
$ws2.Range("F19").Value = "F29"
$ws2.Range("G19").Value = "G29"
$ws2.Range("H19").Value = "H29"
$ws2.Range("B20").Value = "A30"
$ws2.Range("C20").Value = "C30"
$ws2.Range("D20").Value = "D30"
$ws2.Range("E20").Value = "E30"
$ws2.Range("F20").Value = "F30"
$ws2.Range("G20").Value = "G30"
$ws2.Range("H20").Value = "H30"
$ws2.Range("B21").Value = "A31"
$ws2.Range("C21").Value = "C31"
$ws2.Range("D21").Value = "D31"
$ws2.Range("E21").Value = "E31"
$ws2.Range("F21").Value = "F31"
$ws2.Range("G21").Value = "G31"
$ws2.Range("H21").Value = "H31"
$ws2.Range("B22").Value = "A32"
$ws2.Range("C22").Value = "C32"
$ws2.Range("D22").Value = "D32"
$ws2.Range("E22").Value = "E32"
$ws2.Range("F22").Value = "F32"
$ws2.Range("G22").Value = "G32"
$ws2.Range("H22").Value = "H32"
$ws2.Range("B23").Value = "A33"
$ws2.Range("C23").Value = "C33"
$ws2.Range("D23").Value = "D33"
$ws2.Range("E23").Value = "E33"
$ws2.Range("F23").Value = "F33"
$ws2.Range("G23").Value = "G33"
$ws2.Range("H23").Value = "H33"

# Fix up selections: Sheet1 gets a full-range selection (A1:H23) and is no
# longer the active tab; Sheet2 becomes the active tab with H3 selected.
$null = $ws1.Range("A1:H23").Select()
$null = $ws2.Range("H3").Select()
$ws2.Activate()
